# Weekly data refresh: insert two new price observations at the top of the
# Albahaca / Vega Modelo de Temuco daily logic block (rows 409-410), pushing
# all existing rows down by two (old 409..443 become 411..445).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 409; this shifts the existing
# 409:443 block down to 411:445 and carries the D-column date number format
# (style index 2) onto the new rows, matching the rest of the column.
$ws.Rows("409:410").Insert()

# --- New row 409 -----------------------------------------------------
$ws.Cells.Item(409, 1).Value  = 10
$ws.Cells.Item(409, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(409, 3).Value  = "La Araucanía"
$ws.Cells.Item(409, 4).Value  = 45223
$ws.Cells.Item(409, 5).Value  = 9
$ws.Cells.Item(409, 6).Value  = 100112052
$ws.Cells.Item(409, 7).Value  = "Albahaca"
$ws.Cells.Item(409, 8).Value  = "Sin especificar"
$ws.Cells.Item(409, 9).Value  = "Primera"
$ws.Cells.Item(409, 10).Value = 110
$ws.Cells.Item(409, 11).Value = 5000
$ws.Cells.Item(409, 12).Value = 5000
$ws.Cells.Item(409, 13).Value = 5000
$ws.Cells.Item(409, 14).Value = "`$/paquete"
$ws.Cells.Item(409, 15).Value = "Región Metropolitana"
$ws.Cells.Item(409, 16).Value = 5000
$ws.Cells.Item(409, 17).Value = 1
$ws.Cells.Item(409, 18).Value = "Hortaliza"

# --- New row 410 -----------------------------------------------------
$ws.Cells.Item(410, 1).Value  = 10
$ws.Cells.Item(410, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(410, 3).Value  = "La Araucanía"
$ws.Cells.Item(410, 4).Value  = 45223
$ws.Cells.Item(410, 5).Value  = 9
$ws.Cells.Item(410, 6).Value  = 100112052
$ws.Cells.Item(410, 7).Value  = "Albahaca"
$ws.Cells.Item(410, 8).Value  = "Sin especificar"
$ws.Cells.Item(410, 9).Value  = "Primera"
$ws.Cells.Item(410, 10).Value = 60
$ws.Cells.Item(410, 11).Value = 6000
$ws.Cells.Item(410, 12).Value = 6000
$ws.Cells.Item(410, 13).Value = 6000
$ws.Cells.Item(410, 14).Value = "`$/paquete"
$ws.Cells.Item(410, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(410, 16).Value = 6000
$ws.Cells.Item(410, 17).Value = 1
$ws.Cells.Item(410, 18).Value = "Hortaliza"

Write-Output "rows inserted and populated"
